$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$col = $ws.Columns.Item(4)
$col.Insert()

$eRange = $ws.Range("E5:E102")
$eRange.Copy()
$dRange = $ws.Range("D5:D102")
$dRange.PasteSpecial(-4122)
Write-Host "ok"
